# OperationScenario_Component_Building.xlsx
# "behavior results inserted into operation model"
#
# The single data sheet holds three building rows (A2:P4). The edit:
#   - overwrites the per-row "person_num" (col E) with 1 / 2 / 4
#   - overwrites "Hop" (col G) and "Htr_w" (col H) on every row with the
#     newly computed behavior-model results (286.528 / 115.621)
#   - drops the leftover direct cell formatting (a no-op "applyFill" style)
#     that was sitting on A2:P4 and the blank trailing row 5
#   - removes the now-empty trailing row 5
#   - shrinks the AutoFilter (and its backing _FilterDatabase defined name)
#     from A1:P4 down to A1:P2
#   - leaves the view zoomed to 200% with B9 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the AutoFilter at A1:P2 -----------------------------------
# Simply re-applying AutoFilter over "A1:P2" auto-grows back to the
# contiguous data block (A1:P4) because rows 3:4 are still full. Toggle the
# existing filter off, blank out the rows below the new extent so nothing is
# adjacent, then turn the filter back on over the smaller range.
[void]$ws.AutoFilter.Range().AutoFilter()
$ws.Range("A3:P5").ClearContents()
[void]$ws.Range("A1:P2").AutoFilter()

# --- 2. Write the updated behavior-model results ---------------------------
# Row 2 keeps its original A-D values; only E/G/H change.
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 286.528
$ws.Range("H2").Value = 115.621

# Rows 3 and 4 were blanked out in step 1 (so the AutoFilter resize wouldn't
# re-absorb them) - re-enter them in full, cell by cell, with the refreshed
# numbers.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "SFH"
$ws.Range("C3").Value = 1949
$ws.Range("D3").Value = 1957
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 174.14
$ws.Range("G3").Value = 286.528
$ws.Range("H3").Value = 115.621
$ws.Range("I3").Value = 49.261
$ws.Range("J3").Value = 213505.516
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 10.037
$ws.Range("N3").Value = 4.533
$ws.Range("O3").Value = 1.619
$ws.Range("P3").Value = 21000

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "SFH"
$ws.Range("C4").Value = 1949
$ws.Range("D4").Value = 1957
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 174.14
$ws.Range("G4").Value = 286.528
$ws.Range("H4").Value = 115.621
$ws.Range("I4").Value = 49.261
$ws.Range("J4").Value = 213505.516
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 10.037
$ws.Range("N4").Value = 4.533
$ws.Range("O4").Value = 1.619
$ws.Range("P4").Value = 21000

# --- 3. Drop the leftover direct formatting ---------------------------------
$ws.Rows("2:4").ClearFormats()
$ws.Rows("5").Delete()

# --- 4. Keep the _FilterDatabase defined name in sync with the AutoFilter --
$wb.Names.Item(1).RefersTo = "=OperationScenario_Component_Bui!`$A`$1:`$P`$2"

# --- 5. View state: zoom to 200% and select B9 ------------------------------
$ws.Range("B9").Select()
$excel.ActiveWindow.Zoom = 200
